$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 4 (pushes old rows 4..29 down to 6..31,
# carrying their existing A/B/C..T content with them unchanged)
$ws.Rows.Item(4).Resize(2).Insert()

# Match the formatting (bold, centered, thin border) that column A already
# uses on every other data row, e.g. row 3 / row 6.
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)

# ---- New row 4: Holden ----
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 1.254701697683301
$ws.Range("D4").Value = 0.7773705505138532
$ws.Range("E4").Value = 0.8704374273154003
$ws.Range("F4").Value = 1.254701697683301
$ws.Range("G4").Value = 0.8646450015213883
$ws.Range("H4").Value = 0.7773705505138532
$ws.Range("I4").Value = 1.086051695758476
$ws.Range("J4").Value = 0.8237832698567501
$ws.Range("K4").Value = 0.7773705505138532
$ws.Range("L4").Value = 0.8704374273154003
$ws.Range("M4").Value = 1.062569562499351
$ws.Range("N4").Value = 1.062569562499351
$ws.Range("O4").Value = 1.070396940252393
$ws.Range("P4").Value = 0.9675032251708516
$ws.Range("Q4").Value = 0.9675032251708516
$ws.Range("R4").Value = 0.919970056506602
$ws.Range("S4").Value = 0.919970056506602
$ws.Range("T4").Value = 0.9461649404415283

# ---- New row 5: Rizzie Spiral ----
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 0.001351506250955896
$ws.Range("D5").Value = 4.492740714395812
$ws.Range("E5").Value = 0.001944713813196694
$ws.Range("F5").Value = 0.001351506250955896
$ws.Range("G5").Value = 1.241745914772147
$ws.Range("H5").Value = 4.492740714395812
$ws.Range("I5").Value = 1.802594376952982
$ws.Range("J5").Value = 0.01102323034015119
$ws.Range("K5").Value = 4.492740714395812
$ws.Range("L5").Value = 0.001944713813196694
$ws.Range("M5").Value = 0.001648110032076296
$ws.Range("N5").Value = 0.001648110032076296
$ws.Range("O5").Value = 0.6019635323390449
$ws.Range("P5").Value = 1.498678978153322
$ws.Range("Q5").Value = 1.498678978153322
$ws.Range("R5").Value = 2.247194412213944
$ws.Range("S5").Value = 2.247194412213944
$ws.Range("T5").Value = 1.258566742754207

# ---- Fix up A/B (the sequential index + row label columns) for every row
#      from 6 down to 31: the C:T simulation values already shifted down by
#      2 rows via the Insert above, but the A (index) / B (label) columns
#      must show the correct sequential position, not the old pushed-down
#      values.
$labels = @(
  "RotRing OmegaMax-90",
  "Equal Angle",
  "Tilt Rotate",
  "CLR",
  "Rizzie Hex",
  "Matthies Hex",
  "Tilt Rotate_Partial",
  "RotRing OmegaMax-60",
  "Equal Angle_Partial",
  "Rizzie Hex_Partial",
  "ND Single",
  "RD Single",
  "TD Single",
  "Morris Single",
  "Ring Perpendicular to ND",
  "Ring Perpendicular to RD",
  "Ring Perpendicular to TD",
  "OffsetFTD",
  "OffsetATD",
  "OffsetF45",
  "OffsetA45",
  "OffsetFRD",
  "OffsetARD",
  "Gaussian Quadrature",
  "Michael-CCHex",
  "Michael-SNHex"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = 6 + $i
    $ws.Range("A$r").Value = 4 + $i
    $ws.Range("B$r").Value = $labels[$i]
}
